$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("M2").Value = 8.424749
$ws.Range("N2").Value = 16.849498
$ws.Range("O2").Value = 0.1980372819757593
$ws.Range("P2").Value = 0.1585662899160533
$ws.Range("Q2").Value = 0.2357272852696667
$ws.Range("R2").Value = 1.414363711618
$ws.Range("S2").Value = 0.1980372819757593
$ws.Range("T2").Value = 0.1585662899160533

# Row 3
$ws.Range("O3").Value = 0.2481739331426511
$ws.Range("P3").Value = 0.2980652388254983
$ws.Range("S3").Value = 0.2481739331426511
$ws.Range("T3").Value = 0.2980652388254983

# Row 4
$ws.Range("M4").Value = 5.743874333333333
$ws.Range("N4").Value = 17.231623
$ws.Range("O4").Value = 0.1350190089916815
$ws.Range("P4").Value = 0.1621623699615343
$ws.Range("Q4").Value = 0.1607155184714444
$ws.Range("R4").Value = 1.446439666243
$ws.Range("S4").Value = 0.1350190089916815
$ws.Range("T4").Value = 0.1621623699615343

# Row 5
$ws.Range("M5").Value = 12.9373935
$ws.Range("N5").Value = 25.874787
$ws.Range("O5").Value = 0.3041142524947457
$ws.Range("P5").Value = 0.2435009622813764
$ws.Range("Q5").Value = 0.3619925825945
$ws.Range("R5").Value = 2.171955495567
$ws.Range("S5").Value = 0.3041142524947457
$ws.Range("T5").Value = 0.2435009622813764

# Row 6
$ws.Range("M6").Value = 0.882742
$ws.Range("N6").Value = 2.648226
$ws.Range("O6").Value = 0.02075027117909931
$ws.Range("P6").Value = 0.02492177343676531
$ws.Range("Q6").Value = 0.02469941540733333
$ws.Range("R6").Value = 0.222294738666
$ws.Range("S6").Value = 0.02075027117909931
$ws.Range("T6").Value = 0.02492177343676531

# Row 7
$ws.Range("M7").Value = 3.994844666666667
$ws.Range("N7").Value = 11.984534
$ws.Range("O7").Value = 0.09390525221606305
$ws.Range("P7").Value = 0.1127833655787726
$ws.Range("Q7").Value = 0.1117770853882222
$ws.Range("R7").Value = 1.005993768494
$ws.Range("S7").Value = 0.09390525221606305
$ws.Range("T7").Value = 0.1127833655787726
